$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting for new rows 102:103 from row 101 ---
$ws.Range("A101:AB101").Copy()
$ws.Range("A102:AB103").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Cell value updates ---
$ws.Range("AA11").Value = 1
$ws.Range("AB11").Value = -1
$ws.Range("B11").Value = 6227815
$ws.Range("E11").Value = "HFX Wanderers"
$ws.Range("F11").Value = "Cavalry FC"
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 1
$ws.Range("J11").Value = 2.6
$ws.Range("K11").Value = 3.2
$ws.Range("L11").Value = 2.4
$ws.Range("M11").Value = 3.3
$ws.Range("N11").Value = 3
$ws.Range("O11").Value = 2.15
$ws.Range("P11").Value = 0.25
$ws.Range("Q11").Value = 1.925
$ws.Range("R11").Value = 1.875
$ws.Range("S11").Value = 2.25
$ws.Range("T11").Value = 2
$ws.Range("U11").Value = 1.8
$ws.Range("V11").Value = 2.3
$ws.Range("Y11").Value = 0.925
$ws.Range("Z11").Value = -1
$ws.Range("AA12").Value = -1
$ws.Range("AB12").Value = 0.825
$ws.Range("B12").Value = 6240280
$ws.Range("E12").Value = "Atletico Ottawa"
$ws.Range("F12").Value = "Vancouver FC"
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 1.571
$ws.Range("K12").Value = 3.4
$ws.Range("L12").Value = 5.5
$ws.Range("M12").Value = 1.444
$ws.Range("N12").Value = 3.8
$ws.Range("O12").Value = 6
$ws.Range("P12").Value = -1.25
$ws.Range("Q12").Value = 1.95
$ws.Range("R12").Value = 1.85
$ws.Range("S12").Value = 2.75
$ws.Range("T12").Value = 1.975
$ws.Range("U12").Value = 1.825
$ws.Range("V12").Value = 0.444
$ws.Range("Y12").Value = -0.5
$ws.Range("Z12").Value = 0.425
$ws.Range("AA101").Value = 0.8
$ws.Range("AB101").Value = -1
$ws.Range("B101").Value = 7803365
$ws.Range("G101").Value = 3
$ws.Range("H101").Value = 1
$ws.Range("I101").Value = "H"
$ws.Range("O101").Value = 4
$ws.Range("S101").Value = 2.5
$ws.Range("T101").Value = 1.8
$ws.Range("U101").Value = 2
$ws.Range("V101").Value = 0.6499999999999999
$ws.Range("W101").Value = -1
$ws.Range("X101").Value = -1
$ws.Range("Y101").Value = 0.8500000000000001
$ws.Range("Z101").Value = -1
$ws.Range("A102").Value = 100
$ws.Range("AA102").Value = -0.5
$ws.Range("AB102").Value = 0.5
$ws.Range("B102").Value = 7802939
$ws.Range("C102").Value = "Canada Premier League"
$ws.Range("D102").Value = 45416.83333333334
$ws.Range("E102").Value = "Pacific FC CA"
$ws.Range("F102").Value = "York United FC"
$ws.Range("G102").Value = 2
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = "H"
$ws.Range("J102").Value = 1.909
$ws.Range("K102").Value = 3.6
$ws.Range("L102").Value = 3.2
$ws.Range("M102").Value = 1.8
$ws.Range("N102").Value = 3.5
$ws.Range("O102").Value = 3.75
$ws.Range("P102").Value = -0.5
$ws.Range("Q102").Value = 1.825
$ws.Range("R102").Value = 1.975
$ws.Range("S102").Value = 2.25
$ws.Range("T102").Value = 1.8
$ws.Range("U102").Value = 2
$ws.Range("V102").Value = 0.8
$ws.Range("W102").Value = -1
$ws.Range("X102").Value = -1
$ws.Range("Y102").Value = 0.825
$ws.Range("Z102").Value = -1
$ws.Range("A103").Value = 101
$ws.Range("AA103").Value = -1
$ws.Range("AB103").Value = 0.8500000000000001
$ws.Range("B103").Value = 7802940
$ws.Range("C103").Value = "Canada Premier League"
$ws.Range("D103").Value = 45417.625
$ws.Range("E103").Value = "Atletico Ottawa"
$ws.Range("F103").Value = "Valour FC"
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = "H"
$ws.Range("J103").Value = 2
$ws.Range("K103").Value = 3.25
$ws.Range("L103").Value = 3.25
$ws.Range("M103").Value = 1.615
$ws.Range("N103").Value = 3.6
$ws.Range("O103").Value = 4.75
$ws.Range("P103").Value = -0.75
$ws.Range("Q103").Value = 1.8
$ws.Range("R103").Value = 2
$ws.Range("S103").Value = 2.5
$ws.Range("T103").Value = 1.95
$ws.Range("U103").Value = 1.85
$ws.Range("V103").Value = 0.615
$ws.Range("W103").Value = -1
$ws.Range("X103").Value = -1
$ws.Range("Y103").Value = 0.8
$ws.Range("Z103").Value = -1
